$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confirmation prompt text in C2 (shared string content change)
$ws.Range("C2").Value = "Confirm to cancel Deletion on Instance Termination?"

# Update the current selection to match the saved workbook view state
$ws.Range("C13:C14").Select() | Out-Null
